# Add columns I (I0) and J (IF) to the worksheet, matching columns H's header style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy the formatting used by the existing header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-13 for columns I and J
$values = @{
    2  = @(5, 5)
    3  = @(6, 6)
    4  = @(6, 7)
    5  = @(3, 4)
    6  = @(5, 5)
    7  = @(3, 3)
    8  = @(7, 7)
    9  = @(8, 8)
    10 = @(9, 9)
    11 = @(3, 3)
    12 = @(8, 8)
    13 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
